# Team profile edit: zero-out several "Balls Faced"/column-B figures,
# fix up the last player's name (typed "pakaya" first, then corrected
# to "GOTA"), and leave the selection where the user last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column B updates (set to 0) ---
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0

# --- Player name correction on row 12 ---
# The user typed "pakaya" first, then replaced it with "GOTA".
$ws.Range("A12").Value = "pakaya"
$ws.Range("A12").Value = "GOTA"

# --- Final selection / scroll position left by the user ---
$ws.Range("G8").Select()
